$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) stores numeric-looking text (e.g. "370.82", "1.28") as
# plain strings in the source workbook. Pre-format the column as Text so
# the COM layer does not silently coerce these assignments to numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '50.847.52'
$ws.Range('E2').Value = '  -1.23%  '

# Row 3
$ws.Range('D3').Value = '2.921.07'
$ws.Range('E3').Value = '  -1.71%  '

# Row 4
$ws.Range('E4').Value = '  +0.11%  '

# Row 5
$ws.Range('D5').Value = '370.82'
$ws.Range('E5').Value = '  -1.91%  '

# Row 6
$ws.Range('D6').Value = '99.29'
$ws.Range('E6').Value = '  -5.20%  '

# Row 7
$ws.Range('D7').Value = '0.528'
$ws.Range('E7').Value = '  -2.18%  '

# Row 8
$ws.Range('E8').Value = '  +0.07%  '

# Row 9
$ws.Range('D9').Value = '0.572'
$ws.Range('E9').Value = '  -3.28%  '

# Row 10
$ws.Range('D10').Value = '35.54'
$ws.Range('E10').Value = '  -4.48%  '

# Row 11
$ws.Range('E11').Value = '  -0.97%  '

# Row 12
$ws.Range('D12').Value = '0.0835'
$ws.Range('E12').Value = '  -0.86%  '

# Row 13
$ws.Range('D13').Value = '3.376.55'
$ws.Range('E13').Value = '  -1.67%  '

# Row 14
$ws.Range('D14').Value = '17.75'
$ws.Range('E14').Value = '  -3.62%  '

# Row 15
$ws.Range('D15').Value = '7.34'
$ws.Range('E15').Value = '  -2.79%  '

# Row 16
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.930.78'
$ws.Range('E16').Value = '  -1.32%  '

# Row 17
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').Value = '11.12'
$ws.Range('E17').Value = '  +50.51%  '

# Row 18
$ws.Range('D18').Value = '0.950'
$ws.Range('E18').Value = '  -1.66%  '

# Row 19
$ws.Range('D19').Value = '50.792.34'
$ws.Range('E19').Value = '  -1.20%  '

# Row 20
$ws.Range('D20').Value = '3.10'
$ws.Range('E20').Value = '  -6.90%  '

# Row 21
$ws.Range('D21').Value = '12.13'
$ws.Range('E21').Value = '  -5.94%  '

# Row 22
$ws.Range('D22').Value = '0.0₃0942'
$ws.Range('E22').Value = '  -2.04%  '

# Row 23
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').Value = '68.00'
$ws.Range('E23').Value = '  -1.99%  '

# Row 24
$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').Value = '261.80'
$ws.Range('E24').Value = '  +0.12%  '

# Row 25
$ws.Range('D25').Value = '3.06'
$ws.Range('E25').Value = '  +8.45%  '

# Row 26
$ws.Range('D26').Value = '7.93'
$ws.Range('E26').Value = '  -1.81%  '

# Row 27
$ws.Range('D27').Value = '7.16'
$ws.Range('E27').Value = '  -5.60%  '

# Row 28
$ws.Range('E28').Value = '  +0.12%  '

# Row 29
$ws.Range('E29').Value = '  -5.46%  '

# Row 30
$ws.Range('D30').Value = '25.27'
$ws.Range('E30').Value = '  -2.25%  '

# Row 31
$ws.Range('E31').Value = '  -3.02%  '

# Row 32
$ws.Range('D32').Value = '9.76'
$ws.Range('E32').Value = '  -1.29%  '

# Row 33
$ws.Range('D33').Value = '50.39'
$ws.Range('E33').Value = '  -1.07%  '

# Row 34
$ws.Range('E34').Value = '  -3.02%  '

# Row 35
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').Value = '32.45'
$ws.Range('E35').Value = '  -7.37%  '

# Row 36
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').Value = '0.0434'
$ws.Range('E36').Value = '  -2.43%  '

# Row 37
$ws.Range('E37').Value = '  -0.05%  '

# Row 38
$ws.Range('E38').Value = '  +0.88%  '

# Row 39
$ws.Range('E39').Value = '  -1.66%  '

# Row 40
$ws.Range('D40').Value = '15.97'
$ws.Range('E40').Value = '  -7.06%  '

# Row 41
$ws.Range('E41').Value = '  -4.94%  '

# Row 42
$ws.Range('E42').Value = '  -6.62%  '

# Row 43
$ws.Range('D43').Value = '119.19'
$ws.Range('E43').Value = '  -4.42%  '

# Row 44
$ws.Range('D44').Value = '20.79'
$ws.Range('E44').Value = '  -4.61%  '

# Row 45
$ws.Range('E45').Value = '  -1.67%  '

# Row 46
$ws.Range('E46').Value = '  -6.93%  '

# Row 47
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '2.28'
$ws.Range('E47').Value = '  -3.61%  '

# Row 48
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '3.18'
$ws.Range('E48').Value = '  -1.24%  '

# Row 49
$ws.Range('D49').Value = '1.962.98'
$ws.Range('E49').Value = '  -3.50%  '

# Row 50
$ws.Range('B50').Value = 'BEAM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range('D50').Value = '0.0320'
$ws.Range('E50').Value = '  -6.12%  '

# Row 51
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').Value = '1.28'
$ws.Range('E51').Value = '  -0.02%  '

# Clear the temporary Text number-format marker so the cells fall back
# to the plain (unstyled) cell format used throughout the rest of the sheet.
$priceRange.Style = "Normal"
